$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(13289,12279,12279,12279,11330,11330,10778,10693,10693,9025,8836,8756,8756,8756,8637,8637,8153,8153,8153,8153,8153,8153,8153,8153,7884,7884,7884,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7882,7729,7729,7729,7729,7729,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7721,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7312,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7310,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293,7293)

$startRow = 2
$endRow = $startRow + $values.Length - 1

$data = New-Object 'object[,]' $values.Length,1
for ($i = 0; $i -lt $values.Length; $i++) {
    $data[$i,0] = $values[$i]
}

$range = $ws.Range("C$startRow" + ":C$endRow")
$range.Value = $data
